$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column width change (A:E from ~22.27 -> ~19.54 "characters").
#     The COM ColumnWidth setter in this host snaps to a 1/6-character
#     grid, so we pick the input that lands on the closest reachable
#     value (19.5) to the target raw width (19.54296875). ---
$ws.Columns("A:E").ColumnWidth = 18.666666666666668

# --- Header text update: "Accuracy" -> "Accuracy (Over 4 Folds)" ---
$ws.Range("A1").Value = "Accuracy (Over 4 Folds)"

# --- Remove the old "K-Nearest Neighbors" row (row 3) by deleting it,
#     shifting the "Random Forest" row (old row 4) up into row 3,
#     and shifting everything below up by one row as well. ---
$ws.Rows("3").Delete()

# Row 3 is now the old "Random Forest" row; update its values
# (label stays "Random Forest", numbers change per the k-fold accuracy update).
$ws.Range("A3").Value = "Random Forest"
$ws.Range("B3").Value = 0.95189999999999997
$ws.Range("C3").Value = 0.67025462962962901
$ws.Range("D3").Value = 0.59519999999999995
$ws.Range("E3").Value = 0.50970000000000004

# --- Add the new "Cluster" row (row 9) under the McNemar vs LogReg block ---
$ws.Range("A9").Value = "Cluster"
$ws.Range("A9").Font.Bold = $false
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0.50248165485657903
$ws.Range("E9").Value = 0.43824309778037002
$ws.Range("B9:E9").NumberFormat = "0.000"

# --- Update selection to match the final cursor position ---
$ws.Range("H13").Select()
